$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.734961060295085
$ws.Range("D2").Value = 9.00900354053976
$ws.Range("E2").Value = 13.54306243573362
$ws.Range("F2").Value = 33.35405442274314
$ws.Range("G2").Value = 3.657218924134185
$ws.Range("I2").Value = 20.57168573832833
$ws.Range("J2").Value = 10.06171222349325
$ws.Range("K2").Value = 11.56035979850169
$ws.Range("M2").Value = 15.95132316739785
$ws.Range("N2").Value = 19.10453818014686
$ws.Range("O2").Value = 25.16966065198412
$ws.Range("B3").Value = 5.615614453262408
$ws.Range("D3").Value = 8.991860043657228
$ws.Range("E3").Value = 13.55089476578844
$ws.Range("F3").Value = 33.39026655727726
$ws.Range("G3").Value = 3.659162760256811
$ws.Range("I3").Value = 20.66588915144029
$ws.Range("J3").Value = 10.08380912083142
$ws.Range("K3").Value = 11.28652091386433
$ws.Range("M3").Value = 15.85201922758775
$ws.Range("N3").Value = 19.16017964157055
$ws.Range("O3").Value = 25.22610042560389
$ws.Range("B4").Value = 5.541944811397091
$ws.Range("D4").Value = 8.982667951417769
$ws.Range("E4").Value = 13.55791254144329
$ws.Range("F4").Value = 33.42012360944247
$ws.Range("G4").Value = 3.66042044886761
$ws.Range("I4").Value = 20.72677624335076
$ws.Range("J4").Value = 10.09849345644314
$ws.Range("K4").Value = 11.11661834451524
$ws.Range("M4").Value = 15.7932351071779
$ws.Range("N4").Value = 19.19603514025241
$ws.Range("O4").Value = 25.26600806172086
$ws.Range("B5").Value = 5.511868443381641
$ws.Range("D5").Value = 8.979260391852112
$ws.Range("E5").Value = 13.56132839244654
$ws.Range("F5").Value = 33.43420545127159
$ws.Range("G5").Value = 3.660949152520907
$ws.Range("I5").Value = 20.75235592689284
$ws.Range("J5").Value = 10.10475857201414
$ws.Range("K5").Value = 11.04703716172188
$ws.Range("M5").Value = 15.76984984126718
$ws.Range("N5").Value = 19.21107303480397
$ws.Range("O5").Value = 25.28358947343139
$ws.Range("B6").Value = 5.506872213844341
$ws.Range("D6").Value = 8.978715079282594
$ws.Range("E6").Value = 13.56192919195764
$ws.Range("F6").Value = 33.43665931190371
$ws.Range("G6").Value = 3.661037922397806
$ws.Range("I6").Value = 20.75664983423842
$ws.Range("J6").Value = 10.10581587839631
$ws.Range("K6").Value = 11.03546543090845
$ws.Range("M6").Value = 15.76600169478066
$ws.Range("N6").Value = 19.21359585896554
$ws.Range("O6").Value = 25.28658843008201
$ws.Range("B7").Value = 5.541539357830133
$ws.Range("D7").Value = 8.982620622575654
$ws.Range("E7").Value = 13.55795635662616
$ws.Range("F7").Value = 33.42030577204206
$ws.Range("G7").Value = 3.660427513549131
$ws.Range("I7").Value = 20.72711810911812
$ws.Range("J7").Value = 10.09857681128386
$ws.Range("K7").Value = 11.11568120917691
$ws.Range("M7").Value = 15.79291739354946
$ws.Range("N7").Value = 19.19623621817477
$ws.Range("O7").Value = 25.26623983451982
$ws.Range("B8").Value = 5.693916615825552
$ws.Range("D8").Value = 9.002817335722625
$ws.Range("E8").Value = 13.54530504837895
$ws.Range("F8").Value = 33.36495733717818
$ws.Range("G8").Value = 3.657875871704977
$ws.Range("I8").Value = 20.60353600478241
$ws.Range("J8").Value = 10.06909964108532
$ws.Range("K8").Value = 11.46636451177091
$ws.Range("M8").Value = 15.9166407018757
$ws.Range("N8").Value = 19.12337300707029
$ws.Range("O8").Value = 25.18802943999991
$ws.Range("B9").Value = 5.987823558966094
$ws.Range("D9").Value = 9.052875199834029
$ws.Range("E9").Value = 13.53798668516391
$ws.Range("F9").Value = 33.316959402831
$ws.Range("G9").Value = 3.653378933700667
$ws.Range("I9").Value = 20.38527599903653
$ws.Range("J9").Value = 10.02014198703981
$ws.Range("K9").Value = 12.13592581345781
$ws.Range("M9").Value = 16.17574480710136
$ws.Range("N9").Value = 18.99385527861457
$ws.Range("O9").Value = 25.07643629943126
$ws.Range("B10").Value = 6.198466879009145
$ws.Range("D10").Value = 9.095833188253131
$ws.Range("E10").Value = 13.54321888352424
$ws.Range("F10").Value = 33.31863912899301
$ws.Range("G10").Value = 3.650380790157785
$ws.Range("I10").Value = 20.23948544819868
$ws.Range("J10").Value = 9.989547170862418
$ws.Range("K10").Value = 12.61163091071535
$ws.Range("M10").Value = 16.37492594970504
$ws.Range("N10").Value = 18.90677123842251
$ws.Range("O10").Value = 25.02003355614719
$ws.Range("B11").Value = 6.292697704801102
$ws.Range("D11").Value = 9.116672576575652
$ws.Range("E11").Value = 13.54788777060069
$ws.Range("F11").Value = 33.3274166421518
$ws.Range("G11").Value = 3.64908256546796
$ws.Range("I11").Value = 20.17629941648531
$ws.Range("J11").Value = 9.97679174989867
$ws.Range("K11").Value = 12.82353533419285
$ws.Range("M11").Value = 16.46717927407105
$ws.Range("N11").Value = 18.86889119868406
$ws.Range("O11").Value = 24.99994957811505
$ws.Range("B12").Value = 6.328117201545615
$ws.Range("D12").Value = 9.124746300088441
$ws.Range("E12").Value = 13.54998326137192
$ws.Range("F12").Value = 33.33189051336984
$ws.Range("G12").Value = 3.648600349538444
$ws.Range("I12").Value = 20.15282148146131
$ws.Range("J12").Value = 9.97212842763995
$ws.Range("K12").Value = 12.90305729714956
$ws.Range("M12").Value = 16.50232600211589
$ws.Range("N12").Value = 18.85479530598404
$ws.Range("O12").Value = 24.99314692057136
$ws.Range("B13").Value = 6.320501235091128
$ws.Range("D13").Value = 9.122999441963078
$ws.Range("E13").Value = 13.54951742298232
$ws.Range("F13").Value = 33.33087588376534
$ws.Range("G13").Value = 3.64870378631657
$ws.Range("I13").Value = 20.15785791886138
$ws.Range("J13").Value = 9.973125339716367
$ws.Range("K13").Value = 12.88596404818115
$ws.Range("M13").Value = 16.49474746610189
$ws.Range("N13").Value = 18.85782007468935
$ws.Range("O13").Value = 24.99457627824525
$ws.Range("B14").Value = 6.295617163428452
$ws.Range("D14").Value = 9.117333178423531
$ws.Range("E14").Value = 13.54805361463016
$ws.Range("F14").Value = 33.32776167845525
$ws.Range("G14").Value = 3.64904270527491
$ws.Range("I14").Value = 20.17435887785678
$ws.Range("J14").Value = 9.976404752568993
$ws.Range("K14").Value = 12.8300925260223
$ws.Range("M14").Value = 16.47006669899566
$ws.Range("N14").Value = 18.86772654749273
$ws.Range("O14").Value = 24.99937382618942
$ws.Range("B15").Value = 6.280339620705942
$ws.Range("D15").Value = 9.113886038299194
$ws.Range("E15").Value = 13.54719958623418
$ws.Range("F15").Value = 33.32600382079349
$ws.Range("G15").Value = 3.649251524880361
$ws.Range("I15").Value = 20.18452465122996
$ws.Range("J15").Value = 9.978435212415144
$ws.Range("K15").Value = 12.7957734976982
$ws.Range("M15").Value = 16.45497595584133
$ws.Range("N15").Value = 18.87382687430296
$ws.Range("O15").Value = 25.00241702925421
$ws.Range("B16").Value = 6.192273616409303
$ws.Range("D16").Value = 9.094497033785482
$ws.Range("E16").Value = 13.54295968070217
$ws.Range("F16").Value = 33.31822657136659
$ws.Range("G16").Value = 3.650466949197971
$ws.Range("I16").Value = 20.24367775242617
$ws.Range("J16").Value = 9.990404134128768
$ws.Range("K16").Value = 12.59768542969699
$ws.Range("M16").Value = 16.36892802042712
$ws.Range("N16").Value = 18.90928161141241
$ws.Range("O16").Value = 25.02145836116438
$ws.Range("B17").Value = 6.13781479821048
$ws.Range("D17").Value = 9.082931862422996
$ws.Range("E17").Value = 13.5409437339473
$ws.Range("F17").Value = 33.31550654608345
$ws.Range("G17").Value = 3.651229353127643
$ws.Range("I17").Value = 20.28076806383074
$ws.Range("J17").Value = 9.998044191743739
$ws.Range("K17").Value = 12.47495793374394
$ws.Range("M17").Value = 16.31654461049531
$ws.Range("N17").Value = 18.93147556764297
$ws.Range("O17").Value = 25.03456815549302
$ws.Range("B18").Value = 6.106343626346635
$ws.Range("D18").Value = 9.07640228671994
$ws.Range("E18").Value = 13.53999975080216
$ws.Range("F18").Value = 33.31469651126212
$ws.Range("G18").Value = 3.651674049338206
$ws.Range("I18").Value = 20.30239658466735
$ws.Range("J18").Value = 10.00254796903794
$ws.Range("K18").Value = 12.40394923550659
$ws.Range("M18").Value = 16.28657134339701
$ws.Range("N18").Value = 18.94440430413465
$ws.Range("O18").Value = 25.0426331534874
$ws.Range("B19").Value = 6.095663723041802
$ws.Range("D19").Value = 9.074212635114039
$ws.Range("E19").Value = 13.53971719951291
$ws.Range("F19").Value = 33.31455188095188
$ws.Range("G19").Value = 3.651825678949852
$ws.Range("I19").Value = 20.30977037061373
$ws.Range("J19").Value = 10.004091671239
$ws.Range("K19").Value = 12.37983731590371
$ws.Range("M19").Value = 16.27645048237358
$ws.Range("N19").Value = 18.94880984081851
$ws.Range("O19").Value = 25.04545388112021
$ws.Range("B20").Value = 6.143627598564919
$ws.Range("D20").Value = 9.08415035664051
$ws.Range("E20").Value = 13.54113603897449
$ws.Range("F20").Value = 33.31571802834398
$ws.Range("G20").Value = 3.651147554486953
$ws.Range("I20").Value = 20.2767892006844
$ws.Range("J20").Value = 9.997219572073329
$ws.Range("K20").Value = 12.48806642187673
$ws.Range("M20").Value = 16.32210491284341
$ws.Range("N20").Value = 18.92909608343554
$ws.Range("O20").Value = 25.03311829271601
$ws.Range("B21").Value = 6.302933640085374
$ws.Range("D21").Value = 9.11899258503869
$ws.Range("E21").Value = 13.54847469673626
$ws.Range("F21").Value = 33.32864520757366
$ws.Range("G21").Value = 3.648942901988319
$ws.Range("I21").Value = 20.16949996601165
$ws.Range("J21").Value = 9.975436982563261
$ws.Range("K21").Value = 12.846523500983
$ws.Range("M21").Value = 16.47731046786124
$ws.Range("N21").Value = 18.86481004135707
$ws.Range("O21").Value = 24.99794287682804
$ws.Range("B22").Value = 6.405497969381964
$ws.Range("D22").Value = 9.142824443259505
$ws.Range("E22").Value = 13.55517887396077
$ws.Range("F22").Value = 33.34379536125952
$ws.Range("G22").Value = 3.647556765969779
$ws.Range("I22").Value = 20.10199808881331
$ws.Range("J22").Value = 9.962173376285861
$ws.Range("K22").Value = 13.07655892544654
$ws.Range("M22").Value = 16.5799729109261
$ws.Range("N22").Value = 18.82424312416549
$ws.Range("O22").Value = 24.9796328444426
$ws.Range("B23").Value = 6.350910296544218
$ws.Range("D23").Value = 9.130009371051939
$ws.Range("E23").Value = 13.5514267398362
$ws.Range("F23").Value = 33.33509723850955
$ws.Range("G23").Value = 3.648291580051163
$ws.Range("I23").Value = 20.13778608094348
$ws.Range("J23").Value = 9.969163508498408
$ws.Range("K23").Value = 12.95419567193239
$ws.Range("M23").Value = 16.52507583347837
$ws.Range("N23").Value = 18.84576231902401
$ws.Range("O23").Value = 24.9889767995223
$ws.Range("B24").Value = 6.141000133683564
$ws.Range("D24").Value = 9.083599103189426
$ws.Range("E24").Value = 13.54104842802052
$ws.Range("F24").Value = 33.31562006917329
$ws.Range("G24").Value = 3.651184515775271
$ws.Range("I24").Value = 20.27858709514426
$ws.Range("J24").Value = 9.997592035585468
$ws.Range("K24").Value = 12.48214147344981
$ws.Range("M24").Value = 16.31959065395913
$ws.Range("N24").Value = 18.93017132131316
$ws.Range("O24").Value = 25.03377213105541
$ws.Range("B25").Value = 5.909083918445133
$ws.Range("D25").Value = 9.038233101999319
$ws.Range("E25").Value = 13.53809883303141
$ws.Range("F25").Value = 33.323455837749
$ws.Range("G25").Value = 3.654541549511444
$ws.Range("I25").Value = 20.44175494921521
$ws.Range("J25").Value = 10.03244115673958
$ws.Range("K25").Value = 11.95728735771992
$ws.Range("M25").Value = 16.10401028940759
$ws.Range("N25").Value = 19.02747016693486
$ws.Range("O25").Value = 25.10213960820894
